# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 07:34"

# Row 4 - Estados Unidos (A4 = "Estados Unidos")
$ws.Range("B4").Value = 1292850
$ws.Range("C4").Value = 227
$ws.Range("D4").Value = 217251
$ws.Range("E4").Value = 998661
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 76938

# Row 65 - Hungria
$ws.Range("B65").Value = 3178
$ws.Range("C65").Value = 28
$ws.Range("D65").Value = 865
$ws.Range("E65").Value = 1921
$ws.Range("F65").Value = 74
$ws.Range("G65").Value = 9
$ws.Range("H65").Value = 392

# Row 79 - Bulgaria
$ws.Range("B79").Value = 1865
$ws.Range("C79").Value = 36
$ws.Range("D79").Value = 401
$ws.Range("E79").Value = 1380
